$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H28").Value = 10951.917
$ws.Range("I28").Value = 142.3
$ws.Range("J28").Value = 65000
$ws.Range("K28").Value = 142.3
$ws.Range("L28").Value = 65000
$ws.Range("M28").Value = 342.7
$ws.Range("N28").Value = -65970

$ws.Range("H33").Value = 1692.1852
$ws.Range("I33").Value = 1027.56
$ws.Range("J33").Value = 10000
$ws.Range("K33").Value = 1027.56
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = -798.5599999999999
$ws.Range("N33").Value = -10458

$ws.Range("H43").Value = 1304.5454
$ws.Range("I43").Value = 900
$ws.Range("J43").Value = 1535.7142
$ws.Range("K43").Value = 900
$ws.Range("L43").Value = 1535.7142
$ws.Range("M43").Value = -831
$ws.Range("N43").Value = -1673.7142

$ws.Range("H116").Value = 2909.0908
$ws.Range("I116").Value = 2812.5
$ws.Range("J116").Value = 3166.6667
$ws.Range("K116").Value = 2812.5
$ws.Range("L116").Value = 3166.6667
$ws.Range("M116").Value = 629.5
$ws.Range("N116").Value = -10050.6667

$ws.Range("H138").Value = 5653.207
$ws.Range("I138").Value = 7199.3335
$ws.Range("J138").Value = 5249.8696
$ws.Range("K138").Value = 21598.0005
$ws.Range("L138").Value = 15749.6088
$ws.Range("M138").Value = -16458.0005
$ws.Range("N138").Value = -26029.6088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H44").Value = 19554
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 19554
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 19554
$ws.Range("N44").Value = -20530

$ws.Range("H68").Value = 78033
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 78033
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 78033
$ws.Range("N68").Value = -79655

$ws.Range("H71").Value = 78033
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 78033
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 234099
$ws.Range("N71").Value = -242211

$ws.Range("H132").Value = 3921.7856
$ws.Range("I132").Value = 3037.182
$ws.Range("J132").Value = 7165.3335
$ws.Range("K132").Value = 9111.545999999998
$ws.Range("L132").Value = 21496.0005
$ws.Range("M132").Value = -6581.545999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 70009
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 70009
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 70009
$ws.Range("N14").Value = -70353

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1500.091
$ws.Range("I31").Value = 1555.8422
$ws.Range("J31").Value = 1424.4286
$ws.Range("K31").Value = 1555.8422
$ws.Range("L31").Value = 1424.4286
$ws.Range("M31").Value = -1260.8422
$ws.Range("N31").Value = -2014.4286

$ws.Range("H34").Value = 1500.091
$ws.Range("I34").Value = 1555.8422
$ws.Range("J34").Value = 1424.4286
$ws.Range("K34").Value = 1555.8422
$ws.Range("L34").Value = 1424.4286
$ws.Range("M34").Value = -1353.8422
$ws.Range("N34").Value = -1828.4286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 817.5873
$ws.Range("I68").Value = 585.8182
$ws.Range("J68").Value = 941.95123
$ws.Range("K68").Value = 1757.4546
$ws.Range("L68").Value = 2825.85369
$ws.Range("M68").Value = -946.4546
$ws.Range("N68").Value = -4447.85369

$ws.Range("H71").Value = 817.5873
$ws.Range("I71").Value = 585.8182
$ws.Range("J71").Value = 941.95123
$ws.Range("K71").Value = 5272.3638
$ws.Range("L71").Value = 8477.56107
$ws.Range("M71").Value = -1216.3638
$ws.Range("N71").Value = -16589.56107

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 40000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 40000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 40000
$ws.Range("N19").Value = -40576

$ws.Range("H42").Value = 98290
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 98290
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 98290
$ws.Range("N42").Value = -99260

$ws.Range("H52").Value = 14000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 14000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 14000
$ws.Range("N52").Value = -14518

$ws.Range("H53").Value = 8043
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 8043
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 8043
$ws.Range("N53").Value = -9305

$ws.Range("H58").Value = 12608.2
$ws.Range("I58").Value = 3041
$ws.Range("J58").Value = 15000
$ws.Range("K58").Value = 3041
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -2764
$ws.Range("N58").Value = -15554

$ws.Range("H86").Value = 59969.75
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 59969.75
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 59969.75
$ws.Range("N86").Value = -62341.75

$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H89").Value = 59969.75
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 59969.75
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 179909.25
$ws.Range("N89").Value = -191765.25

$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 802.6
$ws.Range("I107").Value = 603.25
$ws.Range("J107").Value = 1600
$ws.Range("K107").Value = 603.25
$ws.Range("L107").Value = 1600
$ws.Range("M107").Value = 1316.75
$ws.Range("N107").Value = -5440

$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H113").Value = 1770.3334
$ws.Range("I113").Value = 1155.5
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1155.5
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1014.5
$ws.Range("N113").Value = -7340

$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H115").Value = 98290
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 98290
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 98290
$ws.Range("N115").Value = -100640

$ws.Range("H117").Value = 100310
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 100310
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 100310
$ws.Range("N117").Value = -107194

$ws.Range("H118").Value = 77155
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 77155
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 77155
$ws.Range("N118").Value = -80469

$ws.Range("H121").Value = 23200
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 23200
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 23200
$ws.Range("N121").Value = -26694

$ws.Range("H122").Value = 3699
$ws.Range("I122").Value = 2497.5
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 7492.5
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -5042.5
$ws.Range("N122").Value = -18400

$ws.Range("H135").Value = 63125
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 63125
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 63125
$ws.Range("N135").Value = -73265

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 38201
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 38201
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 38201
$ws.Range("N3").Value = -38425

$ws.Range("H15").Value = 38201
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 38201
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 38201
$ws.Range("N15").Value = -38541

$ws.Range("H57").Value = 3041
$ws.Range("I57").Value = 3041
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 3041
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -2475

$ws.Range("H69").Value = 100163
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 100163
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 100163
$ws.Range("N69").Value = -101785

$ws.Range("H70").Value = 87500
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 87500
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 87500
$ws.Range("N70").Value = -88040

$ws.Range("H72").Value = 100163
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 100163
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 300489
$ws.Range("N72").Value = -308601

$ws.Range("H73").Value = 87500
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 87500
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 87500
$ws.Range("N73").Value = -89372

$ws.Range("H132").Value = 4026.111
$ws.Range("I132").Value = 3325.5833
$ws.Range("J132").Value = 4586.533
$ws.Range("K132").Value = 9976.749899999999
$ws.Range("L132").Value = 13759.599
$ws.Range("M132").Value = -7446.749899999999
$ws.Range("N132").Value = -18819.599

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 70002
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 70002
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 70002
$ws.Range("N5").Value = -70226

$ws.Range("H14").Value = 3335993
$ws.Range("I14").Value = 3989.5
$ws.Range("J14").Value = 10000000
$ws.Range("K14").Value = 3989.5
$ws.Range("L14").Value = 10000000
$ws.Range("M14").Value = -3821.5
$ws.Range("N14").Value = -10000336

$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws.Range("H81").Value = 5267.143
$ws.Range("I81").Value = 7717.5
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 15435
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -14374
$ws.Range("N81").Value = -6122

$ws.Range("H84").Value = 5267.143
$ws.Range("I84").Value = 7717.5
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 77175
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -71871
$ws.Range("N84").Value = -30608

$ws.Range("H132").Value = 9725104
$ws.Range("I132").Value = 3107.7778
$ws.Range("J132").Value = 24308098
$ws.Range("K132").Value = 9323.3334
$ws.Range("L132").Value = 72924294
$ws.Range("M132").Value = -6793.3334

$ws.Range("H133").Value = 60715
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 60715
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 60715
$ws.Range("N133").Value = -70835
